$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: write G-column (time strings) for all new rows, in row order,
# so new shared-string entries are created in the exact order the source
# workbook used (00:57:01 .. 01:01:26).
$gValues = @(
    "00:57:01",
    "00:58:39",
    "00:58:39",
    "01:01:13",
    "01:00:58",
    "00:39:05",
    "00:56:40",
    "01:00:10",
    "00:38:24",
    "00:56:48",
    "00:58:11",
    "01:01:26"
)
for ($i = 0; $i -lt $gValues.Length; $i++) {
    $ws.Cells.Item(375 + $i, 7).Value = $gValues[$i]
}

# Step 2: write D-column ("J-2") for all new rows — this is the last new
# unique shared string added (comes after all the G-column times).
for ($i = 0; $i -lt 12; $i++) {
    $ws.Cells.Item(375 + $i, 4).Value = "J-2"
}

# Step 3: fill in the remaining columns for each row (A, B, C, E, F, H..V).
# These only reference strings/numbers already present, so order no longer matters.
# Row 375
$ws.Cells.Item(375, 1).Value = "Entrainement"
$ws.Cells.Item(375, 2).Value = 45897
$ws.Cells.Item(375, 3).Value = "Global"
$ws.Cells.Item(375, 5).Value = "Jeremie Laurent"
$ws.Cells.Item(375, 6).Value = "left forward"
$ws.Cells.Item(375, 8).Value = 3.15
$ws.Cells.Item(375, 9).Value = 0.3
$ws.Cells.Item(375, 10).Value = 2.84
$ws.Cells.Item(375, 11).Value = 0.14
$ws.Cells.Item(375, 12).Value = 0.1
$ws.Cells.Item(375, 13).Value = 0.07
$ws.Cells.Item(375, 14).Value = 0
$ws.Cells.Item(375, 15).Value = 4
$ws.Cells.Item(375, 16).Value = 3.2
$ws.Cells.Item(375, 17).Value = 30.04
$ws.Cells.Item(375, 18).Value = 4.58
$ws.Cells.Item(375, 19).Value = 15
$ws.Cells.Item(375, 20).Value = 4
$ws.Cells.Item(375, 21).Value = 5
$ws.Cells.Item(375, 22).Value = 2

# Row 376
$ws.Cells.Item(376, 1).Value = "Entrainement"
$ws.Cells.Item(376, 2).Value = 45897
$ws.Cells.Item(376, 3).Value = "Global"
$ws.Cells.Item(376, 5).Value = "Romain Thunet"
$ws.Cells.Item(376, 6).Value = "center back"
$ws.Cells.Item(376, 8).Value = 2.75
$ws.Cells.Item(376, 9).Value = 0.2
$ws.Cells.Item(376, 10).Value = 2.54
$ws.Cells.Item(376, 11).Value = 0.1
$ws.Cells.Item(376, 12).Value = 0.08
$ws.Cells.Item(376, 13).Value = 0.02
$ws.Cells.Item(376, 14).Value = 0
$ws.Cells.Item(376, 15).Value = 3
$ws.Cells.Item(376, 16).Value = 2.68
$ws.Cells.Item(376, 17).Value = 27.15
$ws.Cells.Item(376, 18).Value = 3.5
$ws.Cells.Item(376, 19).Value = 3
$ws.Cells.Item(376, 20).Value = 0
$ws.Cells.Item(376, 21).Value = 2
$ws.Cells.Item(376, 22).Value = 1

# Row 377
$ws.Cells.Item(377, 1).Value = "Entrainement"
$ws.Cells.Item(377, 2).Value = 45897
$ws.Cells.Item(377, 3).Value = "Global"
$ws.Cells.Item(377, 5).Value = "Yoan Zouma"
$ws.Cells.Item(377, 6).Value = "center back"
$ws.Cells.Item(377, 8).Value = 2.42
$ws.Cells.Item(377, 9).Value = 0.21
$ws.Cells.Item(377, 10).Value = 2.21
$ws.Cells.Item(377, 11).Value = 0.08
$ws.Cells.Item(377, 12).Value = 0.09
$ws.Cells.Item(377, 13).Value = 0.05
$ws.Cells.Item(377, 14).Value = 0
$ws.Cells.Item(377, 15).Value = 3
$ws.Cells.Item(377, 16).Value = 2.25
$ws.Cells.Item(377, 17).Value = 28.56
$ws.Cells.Item(377, 18).Value = 4.31
$ws.Cells.Item(377, 19).Value = 0
$ws.Cells.Item(377, 20).Value = 1
$ws.Cells.Item(377, 21).Value = 0
$ws.Cells.Item(377, 22).Value = 1

# Row 378
$ws.Cells.Item(378, 1).Value = "Entrainement"
$ws.Cells.Item(378, 2).Value = 45897
$ws.Cells.Item(378, 3).Value = "Global"
$ws.Cells.Item(378, 5).Value = "Omar Benyounes"
$ws.Cells.Item(378, 6).Value = "center midfield"
$ws.Cells.Item(378, 8).Value = 3.57
$ws.Cells.Item(378, 9).Value = 0.3
$ws.Cells.Item(378, 10).Value = 3.26
$ws.Cells.Item(378, 11).Value = 0.16
$ws.Cells.Item(378, 12).Value = 0.08
$ws.Cells.Item(378, 13).Value = 0.06
$ws.Cells.Item(378, 14).Value = 0
$ws.Cells.Item(378, 15).Value = 5
$ws.Cells.Item(378, 16).Value = 3.4
$ws.Cells.Item(378, 17).Value = 27.27
$ws.Cells.Item(378, 18).Value = 4.74
$ws.Cells.Item(378, 19).Value = 18
$ws.Cells.Item(378, 20).Value = 8
$ws.Cells.Item(378, 21).Value = 4
$ws.Cells.Item(378, 22).Value = 2

# Row 379
$ws.Cells.Item(379, 1).Value = "Entrainement"
$ws.Cells.Item(379, 2).Value = 45897
$ws.Cells.Item(379, 3).Value = "Global"
$ws.Cells.Item(379, 5).Value = "Ilan Ihaddadene"
$ws.Cells.Item(379, 6).Value = "center midfield"
$ws.Cells.Item(379, 8).Value = 3.26
$ws.Cells.Item(379, 9).Value = 0.31
$ws.Cells.Item(379, 10).Value = 2.95
$ws.Cells.Item(379, 11).Value = 0.18
$ws.Cells.Item(379, 12).Value = 0.13
$ws.Cells.Item(379, 13).Value = 0
$ws.Cells.Item(379, 14).Value = 0
$ws.Cells.Item(379, 15).Value = 0
$ws.Cells.Item(379, 16).Value = 3.12
$ws.Cells.Item(379, 17).Value = 24.75
$ws.Cells.Item(379, 18).Value = 5.06
$ws.Cells.Item(379, 19).Value = 14
$ws.Cells.Item(379, 20).Value = 8
$ws.Cells.Item(379, 21).Value = 9
$ws.Cells.Item(379, 22).Value = 1

# Row 380
$ws.Cells.Item(380, 1).Value = "Entrainement"
$ws.Cells.Item(380, 2).Value = 45897
$ws.Cells.Item(380, 3).Value = "Global"
$ws.Cells.Item(380, 5).Value = "Emmanuel Valey"
$ws.Cells.Item(380, 6).Value = "left forward"
$ws.Cells.Item(380, 8).Value = 2.43
$ws.Cells.Item(380, 9).Value = 0.31
$ws.Cells.Item(380, 10).Value = 2.11
$ws.Cells.Item(380, 11).Value = 0.15
$ws.Cells.Item(380, 12).Value = 0.06
$ws.Cells.Item(380, 13).Value = 0.1
$ws.Cells.Item(380, 14).Value = 0
$ws.Cells.Item(380, 15).Value = 5
$ws.Cells.Item(380, 16).Value = 3.61
$ws.Cells.Item(380, 17).Value = 29.93
$ws.Cells.Item(380, 18).Value = 4.84
$ws.Cells.Item(380, 19).Value = 10
$ws.Cells.Item(380, 20).Value = 3
$ws.Cells.Item(380, 21).Value = 4
$ws.Cells.Item(380, 22).Value = 0

# Row 381
$ws.Cells.Item(381, 1).Value = "Entrainement"
$ws.Cells.Item(381, 2).Value = 45897
$ws.Cells.Item(381, 3).Value = "Global"
$ws.Cells.Item(381, 5).Value = "Amir Kherrab"
$ws.Cells.Item(381, 6).Value = "center midfield"
$ws.Cells.Item(381, 8).Value = 3.2
$ws.Cells.Item(381, 9).Value = 0.41
$ws.Cells.Item(381, 10).Value = 2.79
$ws.Cells.Item(381, 11).Value = 0.22
$ws.Cells.Item(381, 12).Value = 0.14
$ws.Cells.Item(381, 13).Value = 0.07
$ws.Cells.Item(381, 14).Value = 0
$ws.Cells.Item(381, 15).Value = 4
$ws.Cells.Item(381, 16).Value = 3.26
$ws.Cells.Item(381, 17).Value = 27.43
$ws.Cells.Item(381, 18).Value = 4.42
$ws.Cells.Item(381, 19).Value = 19
$ws.Cells.Item(381, 20).Value = 5
$ws.Cells.Item(381, 21).Value = 4
$ws.Cells.Item(381, 22).Value = 0

# Row 382
$ws.Cells.Item(382, 1).Value = "Entrainement"
$ws.Cells.Item(382, 2).Value = 45897
$ws.Cells.Item(382, 3).Value = "Global"
$ws.Cells.Item(382, 5).Value = "Fareh Wael"
$ws.Cells.Item(382, 6).Value = "center midfield"
$ws.Cells.Item(382, 8).Value = 3.22
$ws.Cells.Item(382, 9).Value = 0.24
$ws.Cells.Item(382, 10).Value = 2.98
$ws.Cells.Item(382, 11).Value = 0.13
$ws.Cells.Item(382, 12).Value = 0.07
$ws.Cells.Item(382, 13).Value = 0.04
$ws.Cells.Item(382, 14).Value = 0
$ws.Cells.Item(382, 15).Value = 4
$ws.Cells.Item(382, 16).Value = 3.12
$ws.Cells.Item(382, 17).Value = 27.61
$ws.Cells.Item(382, 18).Value = 5.07
$ws.Cells.Item(382, 19).Value = 13
$ws.Cells.Item(382, 20).Value = 8
$ws.Cells.Item(382, 21).Value = 7
$ws.Cells.Item(382, 22).Value = 1

# Row 383
$ws.Cells.Item(383, 1).Value = "Entrainement"
$ws.Cells.Item(383, 2).Value = 45897
$ws.Cells.Item(383, 3).Value = "Global"
$ws.Cells.Item(383, 5).Value = "Amine Taiar"
$ws.Cells.Item(383, 6).Value = "center back"
$ws.Cells.Item(383, 8).Value = 2.39
$ws.Cells.Item(383, 9).Value = 0.18
$ws.Cells.Item(383, 10).Value = 2.21
$ws.Cells.Item(383, 11).Value = 0.08
$ws.Cells.Item(383, 12).Value = 0.08
$ws.Cells.Item(383, 13).Value = 0.03
$ws.Cells.Item(383, 14).Value = 0
$ws.Cells.Item(383, 15).Value = 2
$ws.Cells.Item(383, 16).Value = 3.61
$ws.Cells.Item(383, 17).Value = 28.33
$ws.Cells.Item(383, 18).Value = 3.61
$ws.Cells.Item(383, 19).Value = 5
$ws.Cells.Item(383, 20).Value = 0
$ws.Cells.Item(383, 21).Value = 0
$ws.Cells.Item(383, 22).Value = 0

# Row 384
$ws.Cells.Item(384, 1).Value = "Entrainement"
$ws.Cells.Item(384, 2).Value = 45897
$ws.Cells.Item(384, 3).Value = "Global"
$ws.Cells.Item(384, 5).Value = "Sofiane Belle"
$ws.Cells.Item(384, 6).Value = "left forward"
$ws.Cells.Item(384, 8).Value = 2.85
$ws.Cells.Item(384, 9).Value = 0.27
$ws.Cells.Item(384, 10).Value = 2.57
$ws.Cells.Item(384, 11).Value = 0.09
$ws.Cells.Item(384, 12).Value = 0.07
$ws.Cells.Item(384, 13).Value = 0.11
$ws.Cells.Item(384, 14).Value = 0.01
$ws.Cells.Item(384, 15).Value = 6
$ws.Cells.Item(384, 16).Value = 2.68
$ws.Cells.Item(384, 17).Value = 30.63
$ws.Cells.Item(384, 18).Value = 4.72
$ws.Cells.Item(384, 19).Value = 6
$ws.Cells.Item(384, 20).Value = 1
$ws.Cells.Item(384, 21).Value = 7
$ws.Cells.Item(384, 22).Value = 1

# Row 385
$ws.Cells.Item(385, 1).Value = "Entrainement"
$ws.Cells.Item(385, 2).Value = 45897
$ws.Cells.Item(385, 3).Value = "Global"
$ws.Cells.Item(385, 5).Value = "Naim Ighbane"
$ws.Cells.Item(385, 6).Value = "center back"
$ws.Cells.Item(385, 8).Value = 3.25
$ws.Cells.Item(385, 9).Value = 0.17
$ws.Cells.Item(385, 10).Value = 3.08
$ws.Cells.Item(385, 11).Value = 0.05
$ws.Cells.Item(385, 12).Value = 0.06
$ws.Cells.Item(385, 13).Value = 0.06
$ws.Cells.Item(385, 14).Value = 0
$ws.Cells.Item(385, 15).Value = 5
$ws.Cells.Item(385, 16).Value = 2.3
$ws.Cells.Item(385, 17).Value = 27.79
$ws.Cells.Item(385, 18).Value = 3.2
$ws.Cells.Item(385, 19).Value = 5
$ws.Cells.Item(385, 20).Value = 0
$ws.Cells.Item(385, 21).Value = 6
$ws.Cells.Item(385, 22).Value = 0

# Row 386
$ws.Cells.Item(386, 1).Value = "Entrainement"
$ws.Cells.Item(386, 2).Value = 45897
$ws.Cells.Item(386, 3).Value = "Global"
$ws.Cells.Item(386, 5).Value = "Mattheo Haon"
$ws.Cells.Item(386, 6).Value = "right back"
$ws.Cells.Item(386, 8).Value = 3.84
$ws.Cells.Item(386, 9).Value = 0.46
$ws.Cells.Item(386, 10).Value = 3.37
$ws.Cells.Item(386, 11).Value = 0.25
$ws.Cells.Item(386, 12).Value = 0.14
$ws.Cells.Item(386, 13).Value = 0.07
$ws.Cells.Item(386, 14).Value = 0
$ws.Cells.Item(386, 15).Value = 4
$ws.Cells.Item(386, 16).Value = 3.7
$ws.Cells.Item(386, 17).Value = 29.27
$ws.Cells.Item(386, 18).Value = 4.28
$ws.Cells.Item(386, 19).Value = 26
$ws.Cells.Item(386, 20).Value = 4
$ws.Cells.Item(386, 21).Value = 6
$ws.Cells.Item(386, 22).Value = 0

# Step 4: apply the date number format to column B for the new rows so the
# cell style matches the existing date-formatted column (numFmtId 14).
$ws.Range("B375:B386").NumberFormat = "m/d/yy"
